$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A second "Docentes responsaveis" (professor) entry is being added, and it
# needs to show up ABOVE the existing one (row 13: "519033 - Carlos Yujiro
# Shigue"). Inserting a row at 14 clones row 13's formatting (B/C only, no
# A cell) onto the new row - exactly what both the new and the existing
# entry rows look like in the target layout - so insert there, move the
# existing professor down into it, then overwrite row 13 with the new name.
$ws.Rows(14).Insert()

$ws.Range("B14").Value = $ws.Range("B13").Text
$ws.Range("C14").Value = $ws.Range("C13").Text

$ws.Range("B13").Value = "5840730 - Antonio Jefferson da Silva Machado"
$ws.Range("C13").Value = "5840730 - Antonio Jefferson da Silva Machado"
